$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Purchase_Type" column (column D)
$ws.Columns.Item(4).Insert()

# New column header
$ws.Cells.Item(1, 4).Value = "Phone Number"

# Phone numbers for each participant row (kept as text)
$ws.Range("D2:D4").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "7405802474"
$ws.Cells.Item(3, 4).Value = "7016763640"
$ws.Cells.Item(4, 4).Value = "9429510862"
